$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet previously stored a pandas DataFrame with its integer index
# written out to column A (header row had no entry in A1, data column A
# held 0,1,2,3). The re-exported file drops that index column, so every
# column shifts one to the left: B:E -> A:D.
$ws.Columns.Item(1).Delete()
